$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new weekly extract adds 6 records (Melon / Calameno+Tuna x Extra/Primera/Segunda,
# fecha 44610) at the top of this block. Insert 6 blank rows there; Excel shifts all
# the existing rows (and the sheet dimension) down automatically.
$ws.Range("A610:A615").EntireRow.Insert()

# Fill in the 6 newly inserted rows (610-615) cell by cell.
# row 610: Calameño / Extra
$ws.Cells.Item(610, 1).Value = 9
$ws.Cells.Item(610, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(610, 3).Value = "Metropolitana"
$ws.Cells.Item(610, 4).Value = 44610
$ws.Cells.Item(610, 5).Value = 13
$ws.Cells.Item(610, 6).Value = 100112027
$ws.Cells.Item(610, 7).Value = "Melón"
$ws.Cells.Item(610, 8).Value = "Calameño"
$ws.Cells.Item(610, 9).Value = "Extra"
$ws.Cells.Item(610, 10).Value = 250
$ws.Cells.Item(610, 11).Value = 700
$ws.Cells.Item(610, 12).Value = 800
$ws.Cells.Item(610, 13).Value = 750
$ws.Cells.Item(610, 14).Value = "`$/unidad"
$ws.Cells.Item(610, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(610, 16).Value = 750
$ws.Cells.Item(610, 17).Value = 1
$ws.Cells.Item(610, 18).Value = "Hortaliza"
$ws.Cells.Item(610, 4).NumberFormat = $ws.Cells.Item(616, 4).NumberFormat

# row 611: Calameño / Primera
$ws.Cells.Item(611, 1).Value = 9
$ws.Cells.Item(611, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(611, 3).Value = "Metropolitana"
$ws.Cells.Item(611, 4).Value = 44610
$ws.Cells.Item(611, 5).Value = 13
$ws.Cells.Item(611, 6).Value = 100112027
$ws.Cells.Item(611, 7).Value = "Melón"
$ws.Cells.Item(611, 8).Value = "Calameño"
$ws.Cells.Item(611, 9).Value = "Primera"
$ws.Cells.Item(611, 10).Value = 340
$ws.Cells.Item(611, 11).Value = 500
$ws.Cells.Item(611, 12).Value = 600
$ws.Cells.Item(611, 13).Value = 550
$ws.Cells.Item(611, 14).Value = "`$/unidad"
$ws.Cells.Item(611, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(611, 16).Value = 550
$ws.Cells.Item(611, 17).Value = 1
$ws.Cells.Item(611, 18).Value = "Hortaliza"
$ws.Cells.Item(611, 4).NumberFormat = $ws.Cells.Item(616, 4).NumberFormat

# row 612: Calameño / Segunda
$ws.Cells.Item(612, 1).Value = 9
$ws.Cells.Item(612, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(612, 3).Value = "Metropolitana"
$ws.Cells.Item(612, 4).Value = 44610
$ws.Cells.Item(612, 5).Value = 13
$ws.Cells.Item(612, 6).Value = 100112027
$ws.Cells.Item(612, 7).Value = "Melón"
$ws.Cells.Item(612, 8).Value = "Calameño"
$ws.Cells.Item(612, 9).Value = "Segunda"
$ws.Cells.Item(612, 10).Value = 160
$ws.Cells.Item(612, 11).Value = 400
$ws.Cells.Item(612, 12).Value = 400
$ws.Cells.Item(612, 13).Value = 400
$ws.Cells.Item(612, 14).Value = "`$/unidad"
$ws.Cells.Item(612, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(612, 16).Value = 400
$ws.Cells.Item(612, 17).Value = 1
$ws.Cells.Item(612, 18).Value = "Hortaliza"
$ws.Cells.Item(612, 4).NumberFormat = $ws.Cells.Item(616, 4).NumberFormat

# row 613: Tuna / Extra
$ws.Cells.Item(613, 1).Value = 9
$ws.Cells.Item(613, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(613, 3).Value = "Metropolitana"
$ws.Cells.Item(613, 4).Value = 44610
$ws.Cells.Item(613, 5).Value = 13
$ws.Cells.Item(613, 6).Value = 100112027
$ws.Cells.Item(613, 7).Value = "Melón"
$ws.Cells.Item(613, 8).Value = "Tuna"
$ws.Cells.Item(613, 9).Value = "Extra"
$ws.Cells.Item(613, 10).Value = 340
$ws.Cells.Item(613, 11).Value = 700
$ws.Cells.Item(613, 12).Value = 800
$ws.Cells.Item(613, 13).Value = 750
$ws.Cells.Item(613, 14).Value = "`$/unidad"
$ws.Cells.Item(613, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(613, 16).Value = 750
$ws.Cells.Item(613, 17).Value = 1
$ws.Cells.Item(613, 18).Value = "Hortaliza"
$ws.Cells.Item(613, 4).NumberFormat = $ws.Cells.Item(616, 4).NumberFormat

# row 614: Tuna / Primera
$ws.Cells.Item(614, 1).Value = 9
$ws.Cells.Item(614, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(614, 3).Value = "Metropolitana"
$ws.Cells.Item(614, 4).Value = 44610
$ws.Cells.Item(614, 5).Value = 13
$ws.Cells.Item(614, 6).Value = 100112027
$ws.Cells.Item(614, 7).Value = "Melón"
$ws.Cells.Item(614, 8).Value = "Tuna"
$ws.Cells.Item(614, 9).Value = "Primera"
$ws.Cells.Item(614, 10).Value = 430
$ws.Cells.Item(614, 11).Value = 500
$ws.Cells.Item(614, 12).Value = 600
$ws.Cells.Item(614, 13).Value = 550
$ws.Cells.Item(614, 14).Value = "`$/unidad"
$ws.Cells.Item(614, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(614, 16).Value = 550
$ws.Cells.Item(614, 17).Value = 1
$ws.Cells.Item(614, 18).Value = "Hortaliza"
$ws.Cells.Item(614, 4).NumberFormat = $ws.Cells.Item(616, 4).NumberFormat

# row 615: Tuna / Segunda
$ws.Cells.Item(615, 1).Value = 9
$ws.Cells.Item(615, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(615, 3).Value = "Metropolitana"
$ws.Cells.Item(615, 4).Value = 44610
$ws.Cells.Item(615, 5).Value = 13
$ws.Cells.Item(615, 6).Value = 100112027
$ws.Cells.Item(615, 7).Value = "Melón"
$ws.Cells.Item(615, 8).Value = "Tuna"
$ws.Cells.Item(615, 9).Value = "Segunda"
$ws.Cells.Item(615, 10).Value = 250
$ws.Cells.Item(615, 11).Value = 300
$ws.Cells.Item(615, 12).Value = 400
$ws.Cells.Item(615, 13).Value = 350
$ws.Cells.Item(615, 14).Value = "`$/unidad"
$ws.Cells.Item(615, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(615, 16).Value = 350
$ws.Cells.Item(615, 17).Value = 1
$ws.Cells.Item(615, 18).Value = "Hortaliza"
$ws.Cells.Item(615, 4).NumberFormat = $ws.Cells.Item(616, 4).NumberFormat

$ws.Range("A1").Select()
